# Atualização RDD 10 - Liga Eliminação
# Ajuste e atualização da pontuação da Rodada 10 na Liga Eliminação.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column K (shared string "Rodada 10"), matching the
# bold/centered/bordered header style used by the other round columns
$ws.Range("J1").Copy($ws.Range("K1"))
$ws.Cells.Item(1, 11).Value = "Rodada 10"

# Rodada 10 scores for each team (rows with no score simply had no
# participation in this round, matching the sparse Rodada 9 column)
$ws.Cells.Item(2, 11).Value = 40.77001953125
$ws.Cells.Item(3, 11).Value = 66.06982421875
$ws.Cells.Item(4, 11).Value = 69.27001953125
$ws.Cells.Item(5, 11).Value = 78.52001953125
$ws.Cells.Item(6, 11).Value = 79.8701171875
$ws.Cells.Item(7, 11).Value = 53.3701171875
$ws.Cells.Item(8, 11).Value = 54.469970703125
$ws.Cells.Item(9, 11).Value = 73.56982421875
$ws.Cells.Item(10, 11).Value = 61.469970703125
$ws.Cells.Item(12, 11).Value = 76.97021484375
$ws.Cells.Item(13, 11).Value = 73.669921875
$ws.Cells.Item(14, 11).Value = 57.169921875
$ws.Cells.Item(15, 11).Value = 60.77001953125
$ws.Cells.Item(17, 11).Value = 68.669921875
$ws.Cells.Item(19, 11).Value = 66.47021484375
$ws.Cells.Item(20, 11).Value = 64.5
$ws.Cells.Item(21, 11).Value = 83.8701171875
$ws.Cells.Item(22, 11).Value = 54.27001953125
$ws.Cells.Item(26, 11).Value = 65.47021484375
$ws.Cells.Item(29, 11).Value = 53.969970703125
$ws.Cells.Item(30, 11).Value = 57.6201171875
$ws.Cells.Item(32, 11).Value = 61.070068359375
$ws.Cells.Item(33, 11).Value = 68.97021484375
